$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in columns D (Price) and E (Volume) are stored as plain text
# strings in the source data (e.g. "58.045.89", "  +1.34%  "), so force the
# number format to Text ("@") before assigning, to stop Excel from silently
# re-interpreting/parsing them as numeric values.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.045.89'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.114.48'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.70'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.45'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.16%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.113.11'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.62%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.30'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.386'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.649.71'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.43'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.68%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '58.153.94'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.112.98'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.12'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.89'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.13'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '339.06'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.506'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.18'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.168'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0917'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.64'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.90%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.25'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.87'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.92%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.42%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '21.02'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '154.21'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.62'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.05'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.39'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0668'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.158.31'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.681'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.89'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '36.83'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.47'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +6.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.310.34'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0258'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.75'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.958'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.00'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.11%  '
